$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix up formatting that needs to change BEFORE we remove the K:P columns ---
# D5:G5 currently use style 16 ; H5:J5 already use the style (7) that D5:G5 should move to.
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122)   # xlPasteFormats

# E6:J6 currently use style 8 (same as D6); M6:P6 already use the style (11) that E6:J6 should move to.
$ws.Range("M6").Copy()
$ws.Range("E6:J6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- 2) Update the year header row (D4:J4): 2007..2013 -> 2015..2021 ---
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2016
$ws.Range("F4").Value = 2017
$ws.Range("G4").Value = 2018
$ws.Range("H4").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("J4").Value = 2021

# --- 3) Update row 5 data values ---
$ws.Range("D5").Value = 2.2197193775563164
$ws.Range("E5").Value = 2.1235271668715399
$ws.Range("F5").Value = 2.7818537161298167
$ws.Range("G5").Value = 6.7272960584548969
$ws.Range("H5").Value = 5.1525830614767187
$ws.Range("I5").Value = 4.4774536255935971
$ws.Range("J5").Value = 4.6024666695867751

# --- 4) Update row 6 data values ---
$ws.Range("D6").Value = 2.2322863217945752
$ws.Range("E6").Value = 2.8603553109638966
$ws.Range("F6").Value = 3.113207036164539
$ws.Range("G6").Value = 6.2970593463100784
$ws.Range("H6").Value = 4.8617746111834492
$ws.Range("I6").Value = 2.6715092780025032
$ws.Range("J6").Value = 4.3694509108608912

# --- 5) Remove the now-unused columns K:XFD (old years 2014..2019 + the trailing columns) ---
$ws.Range("K:XFD").Delete()

# --- 6) Narrow the data columns D:J down to the new width (~9.43 chars) ---
for ($c = 4; $c -le 10; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 8.666666
}

# --- 7) Update the remembered selection ---
$ws.Range("K16").Select()
